# "fixed lr at 0.5 for base"
#
# The Learning Rate (column E) for the "base" row and the perc1-3 /
# epoch1-3 / layer1-3 comparison rows (rows 2-11) was pinned to 0.5
# (previously a near-zero placeholder of 0.0005). Rows 9-11 (layer1-3)
# only ever had their learning rate column touched. The Accuracy (F)
# and F1 (G) results for rows 2-8 and 12-14 were re-measured for the
# new learning rate and the literal result values updated accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- row 2: base ---
$ws.Range("E2").Value = 0.5
$ws.Range("F2").Value = 0.7678125
$ws.Range("G2").Value = 0.7688129035678574

# --- row 3: perc1 ---
$ws.Range("E3").Value = 0.5
$ws.Range("F3").Value = 0.776875
$ws.Range("G3").Value = 0.7781903995587898

# --- row 4: perc2 ---
$ws.Range("E4").Value = 0.5
$ws.Range("F4").Value = 0.778125
$ws.Range("G4").Value = 0.7794271079322289

# --- row 5: perc3 ---
$ws.Range("E5").Value = 0.5
$ws.Range("F5").Value = 0.579375
$ws.Range("G5").Value = 0.7336762960031658

# --- row 6: epoch1 ---
$ws.Range("E6").Value = 0.5
$ws.Range("F6").Value = 0.7603124999999999
$ws.Range("G6").Value = 0.7611860454996487

# --- row 7: epoch2 ---
$ws.Range("E7").Value = 0.5
$ws.Range("F7").Value = 0.7884375
$ws.Range("G7").Value = 0.7897762549501113

# --- row 8: epoch3 ---
$ws.Range("E8").Value = 0.5
$ws.Range("F8").Value = 0.7865625000000001
$ws.Range("G8").Value = 0.7877376837176459

# --- row 9: layer1 (learning rate only) ---
$ws.Range("E9").Value = 0.5

# --- row 10: layer2 (learning rate only) ---
$ws.Range("E10").Value = 0.5

# --- row 11: layer3 (learning rate only) ---
$ws.Range("E11").Value = 0.5

# --- row 12: lr1 (Accuracy/F1 only, learning rate itself unchanged at 0.1) ---
$ws.Range("F12").Value = 0.720625
$ws.Range("G12").Value = 0.70841279346272

# --- row 13: lr2 (Accuracy/F1 only, learning rate itself unchanged at 0.5) ---
$ws.Range("F13").Value = 0.7940625
$ws.Range("G13").Value = 0.7945738143674341

# --- row 14: lr3 (Accuracy/F1 only, learning rate itself unchanged at 1) ---
$ws.Range("F14").Value = 0.7925
$ws.Range("G14").Value = 0.7936019512183544

# --- move the active selection cursor to H12, matching the saved view state ---
$null = $ws.Range("H12").Select()
